$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 29 content (A29:O29) so it can be rebuilt across the
# new rows 29-36 below.
$ws.Range("A29:O29").ClearContents()

# --- New example clusters inserted ahead of the existing "guzar*" cluster ---
$ws.Range("A29").Value = "zarurat"
$ws.Range("B29").Value = "zaruratmand"

$ws.Range("A30").Value = "jatila"
$ws.Range("B30").Value = "jatilata"

$ws.Range("A31").Value = "itihas"
$ws.Range("B31").Value = "itihasik"

$ws.Range("A32").Value = "pravahita"
$ws.Range("B32").Value = "pravaha"

$ws.Range("A33").Value = "visheshta"
$ws.Range("B33").Value = "vishesh"

$ws.Range("A34").Value = "adhunik"
$ws.Range("B34").Value = "adhunikata"

# --- Remainder of the "guzar*" cluster (minus "guzarta", now on row 35) ---
$ws.Range("A35").Value = "guzarne"
$ws.Range("B35").Value = "guzarna"
$ws.Range("C35").Value = "guzarin"
$ws.Range("D35").Value = "guzari"
$ws.Range("E35").Value = "guzaren"
$ws.Range("F35").Value = "guzarenge"
$ws.Range("G35").Value = "guzarenga"
$ws.Range("H35").Value = "guzare"
$ws.Range("I35").Value = "guzarne"
$ws.Range("J35").Value = "guzara"
$ws.Range("K35").Value = "guzar"
$ws.Range("L35").Value = "goojarna"

# --- "guzarate"/"guzarata" pair moved down onto their own row 36 ---
$ws.Range("A36").Value = "guzarate"
$ws.Range("B36").Value = "guzarata"

# Update the view: scrolled down toward the bottom of the new data and
# zoomed in, with the active selection on the last populated row.
[void]$ws.Range("A35").Select()
$excel.ActiveWindow.Zoom = 140
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
